# Commit: "add new row ibrahim"
# - Append a new row (A5/B5) with "ibrahim" / "29age"
# - Workbook default font switched from Calibri to Arial (via the Normal
#   cell style, so it re-colors every cell without per-cell overrides)
# - Selection left sitting just past the new data, on A6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-theme the workbook's "Normal" style font (affects all cells that don't
# already carry an explicit font override - i.e. every cell in this sheet).
$normal = $wb.Styles("Normal")
$normal.Font.Name = "Arial"

# New row of data.
$ws.Range("A5").Value = "ibrahim"
$ws.Range("B5").Value = "29age"

# Matches the post-edit selection recorded in the saved workbook.
$ws.Range("A6").Select() | Out-Null
